# Cap Ret per Unit Net Loss.xlsx — apply the eps-us-analysis update
#  * B2 (base retirement-fraction parameter) drops from 0.03 to 0.015
#  * Every other "0.03" cell in column B becomes a formula that refers back
#    to B2 (=$B$2) instead of a hard-coded literal, so they track B2
#  * The "0" cells (plant types excluded from retirement) keep their literal
#    0 value but lose the grey highlight fill they previously had
#  * The grey highlight fill is also stripped from the remaining "now a
#    formula" cells that had it (rows 19-25)
#  * The active sheet/tab moves from "About" to "CRpUNL", with a new
#    selection on CRpUNL (C14) while About keeps its old selection (A18)

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("CRpUNL")

# --- B2: new base value -----------------------------------------------
$ws.Range("B2").Value = 0.015

# --- rows whose literal 0.03 becomes "=$B$2" ---------------------------
$formulaRows = @(3, 4, 5, 7, 8, 13, 14, 15, 19, 20, 21, 22, 23, 24, 25)
foreach ($r in $formulaRows) {
    $ws.Range("B$r").Formula = '=$B$2'
}

# --- every cell in B6:B25 that had the grey "s=4" highlight fill loses it
$unhighlightRows = @(6, 9, 10, 11, 12, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25)
foreach ($r in $unhighlightRows) {
    $ws.Range("B$r").Style = "Normal"
}

# --- selection / active tab --------------------------------------------
$wsAbout.Range("A18").Select()
$ws.Activate()
$ws.Range("C14").Select()
